# Auto-generated Excel COM-interop script
# Applies updated market-price / profit values to the Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""
# Row 15
$ws.Range("H15").Value = 4446.2905
$ws.Range("I15").Value = 4446.2905
$ws.Range("K15").Value = 13338.8715
$ws.Range("M15").Value = -13169.8715
# Row 39
$ws.Range("H39").Value = 732.05
$ws.Range("I39").Value = 825.3
$ws.Range("K39").Value = 2475.9
$ws.Range("M39").Value = -2179.9
# Row 40
$ws.Range("H40").Value = 2011
$ws.Range("I40").Value = 1616.5
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 1616.5
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -1441.5
$ws.Range("N40").Value = -3150
# Row 43
$ws.Range("H43").Value = 9259959
$ws.Range("I43").Value = 749.75
$ws.Range("J43").Value = 27778378
$ws.Range("K43").Value = 749.75
$ws.Range("L43").Value = 27778378
$ws.Range("M43").Value = -680.75
$ws.Range("N43").Value = -27778516
# Row 69
$ws.Range("H69").Value = 3301
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126
# Row 72
$ws.Range("H72").Value = 3301
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632
# Row 76
$ws.Range("H76").Value = 4611
$ws.Range("I76").Value = 5333
$ws.Range("K76").Value = 5333
$ws.Range("M76").Value = -5018
# Row 79
$ws.Range("H79").Value = 4611
$ws.Range("I79").Value = 5333
$ws.Range("K79").Value = 5333
$ws.Range("M79").Value = -4241
# Row 87
$ws.Range("H87").Value = 44797.8
$ws.Range("J87").Value = 44797.8
$ws.Range("L87").Value = 44797.8
$ws.Range("N87").Value = -47293.8
# Row 90
$ws.Range("H90").Value = 44797.8
$ws.Range("J90").Value = 44797.8
$ws.Range("L90").Value = 134393.4
$ws.Range("N90").Value = -146873.4

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1882.125
$ws.Range("I45").Value = 1793.8572
$ws.Range("K45").Value = 1793.8572
$ws.Range("M45").Value = -1416.8572
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
# Row 102
$ws.Range("H102").Value = 27795464
$ws.Range("I102").Value = 33354236
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 33354236
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = -33352614
$ws.Range("N102").Value = -4844
# Row 122
$ws.Range("H122").Value = 1182.8334
$ws.Range("I122").Value = 1019.4
$ws.Range("K122").Value = 3058.2
$ws.Range("M122").Value = -608.1999999999998
# Row 132
$ws.Range("H132").Value = 2298.7556
$ws.Range("I132").Value = 2208.9355
$ws.Range("K132").Value = 6626.806500000001
$ws.Range("M132").Value = -4096.806500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -23992

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1223.8372
$ws.Range("I31").Value = 908.5789
$ws.Range("J31").Value = 1473.4166
$ws.Range("K31").Value = 908.5789
$ws.Range("L31").Value = 1473.4166
$ws.Range("M31").Value = -613.5789
$ws.Range("N31").Value = -2063.4166
# Row 34
$ws.Range("H34").Value = 1223.8372
$ws.Range("I34").Value = 908.5789
$ws.Range("J34").Value = 1473.4166
$ws.Range("K34").Value = 908.5789
$ws.Range("L34").Value = 1473.4166
$ws.Range("M34").Value = -706.5789
$ws.Range("N34").Value = -1877.4166
# Row 41
$ws.Range("H41").Value = 5919.5557
$ws.Range("J41").Value = 27000
$ws.Range("L41").Value = 27000
$ws.Range("N41").Value = -27856
# Row 50
$ws.Range("H50").Value = 27999.666
$ws.Range("J50").Value = 27999.666
$ws.Range("L50").Value = 27999.666
$ws.Range("N50").Value = -29249.666
# Row 62
$ws.Range("H62").Value = 4764395.5
$ws.Range("I62").Value = 2551.3416
$ws.Range("K62").Value = 2551.3416
$ws.Range("M62").Value = -1927.3416
# Row 65
$ws.Range("H65").Value = 4764395.5
$ws.Range("I65").Value = 2551.3416
$ws.Range("K65").Value = 12756.708
$ws.Range("M65").Value = -9636.708000000001
# Row 109
$ws.Range("H109").Value = 8146.5
$ws.Range("J109").Value = 8146.5
$ws.Range("L109").Value = 8146.5
$ws.Range("N109").Value = -10226.5

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("J86").Value = 499.85715
$ws.Range("L86").Value = 1499.57145
$ws.Range("N86").Value = -3871.57145
# Row 89
$ws.Range("J89").Value = 499.85715
$ws.Range("L89").Value = 4498.71435
$ws.Range("N89").Value = -16354.71435
# Row 92
$ws.Range("H92").Value = 1250.5
$ws.Range("J92").Value = 1250.5
$ws.Range("L92").Value = 3751.5
$ws.Range("N92").Value = -6247.5
# Row 98
$ws.Range("H98").Value = 824.1
$ws.Range("J98").Value = 1594.75
$ws.Range("L98").Value = 4784.25
$ws.Range("N98").Value = -7780.25
# Row 137
$ws.Range("H137").Value = 4327.3335
$ws.Range("I137").Value = 870
$ws.Range("J137").Value = 7093.2
$ws.Range("K137").Value = 2610
$ws.Range("L137").Value = 21279.6
$ws.Range("M137").Value = 2490
$ws.Range("N137").Value = -31479.6

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 32147872
$ws.Range("I70").Value = 31254700
$ws.Range("J70").Value = 33338766
$ws.Range("K70").Value = 31254700
$ws.Range("L70").Value = 33338766
$ws.Range("M70").Value = -31254430
$ws.Range("N70").Value = -33339306
# Row 73
$ws.Range("H73").Value = 32147872
$ws.Range("I73").Value = 31254700
$ws.Range("J73").Value = 33338766
$ws.Range("K73").Value = 31254700
$ws.Range("L73").Value = 33338766
$ws.Range("M73").Value = -31253764
$ws.Range("N73").Value = -33340638
# Row 80
$ws.Range("H80").Value = 5550
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5550
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5550
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -7546
# Row 83
$ws.Range("H83").Value = 5550
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5550
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 27750
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -37734
# Row 86
$ws.Range("H86").Value = 32533.5
$ws.Range("J86").Value = 32533.5
$ws.Range("L86").Value = 32533.5
$ws.Range("N86").Value = -34905.5
# Row 89
$ws.Range("H89").Value = 32533.5
$ws.Range("J89").Value = 32533.5
$ws.Range("L89").Value = 97600.5
$ws.Range("N89").Value = -109456.5
# Row 122
$ws.Range("H122").Value = 2884.3809
$ws.Range("I122").Value = 1757.2
$ws.Range("K122").Value = 5271.6
$ws.Range("M122").Value = -2821.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 880.6923
$ws.Range("I22").Value = 408.33334
$ws.Range("J22").Value = 1285.5714
$ws.Range("K22").Value = 408.33334
$ws.Range("L22").Value = 1285.5714
$ws.Range("M22").Value = -113.33334
$ws.Range("N22").Value = -1875.5714
# Row 27
$ws.Range("H27").Value = 880.6923
$ws.Range("I27").Value = 408.33334
$ws.Range("J27").Value = 1285.5714
$ws.Range("K27").Value = 408.33334
$ws.Range("L27").Value = 1285.5714
$ws.Range("M27").Value = -301.33334
$ws.Range("N27").Value = -1499.5714
# Row 46
$ws.Range("H46").Value = 5690.8125
$ws.Range("I46").Value = 418.875
$ws.Range("K46").Value = 418.875
$ws.Range("M46").Value = -230.875
# Row 50
$ws.Range("H50").Value = 8389.333000000001
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 10084
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 10084
$ws.Range("M50").Value = -4363
$ws.Range("N50").Value = -11358
# Row 132
$ws.Range("H132").Value = 65075.375
$ws.Range("I132").Value = 2866.3333
$ws.Range("J132").Value = 102400.8
$ws.Range("K132").Value = 8598.999899999999
$ws.Range("L132").Value = 307202.4
$ws.Range("M132").Value = -6068.999899999999
$ws.Range("N132").Value = -312262.4

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2304.4546
$ws.Range("I81").Value = 1534.9
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 3069.8
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -2008.8
$ws.Range("N81").Value = -22122
# Row 84
$ws.Range("H84").Value = 2304.4546
$ws.Range("I84").Value = 1534.9
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 15349
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -10045
$ws.Range("N84").Value = -110608
# Row 122
$ws.Range("H122").Value = 11819754
$ws.Range("I122").Value = 13001529
$ws.Range("J122").Value = 2002.5
$ws.Range("K122").Value = 39004587
$ws.Range("L122").Value = 6007.5
$ws.Range("M122").Value = -39002137
$ws.Range("N122").Value = -10907.5
